# Resume edit: "Java (front and back end), Android ..." ->
# "Java (proficient in all paradigms), Android ..."
# (commit message: "Added 'proficient in all paradigms'")

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$found = $find.Execute("front and back end", $true, $false, $false, $false, $false, `
               $true, 1, $false, "proficient in all paradigms", 2)

if (-not $found) {
    throw "Could not find target text 'front and back end' to replace."
}

Write-Output "Replaced 'front and back end' with 'proficient in all paradigms': $found"
